$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (old N/O/P shift right to O/P/Q)
$mWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $mWidth

$ws.Range("R8").Select() | Out-Null
